{"js": "// Replace each two-digit-division answer with its new value.\n// A few source strings repeat (e.g. \"97\\u00f79=10, 7\" appears twice\n// with different targets), so replacements are keyed on the list of\n// new values for each distinct source text, applied in the order the\n// matches occur in the document (Body.search returns matches in\n// document order).\nconst replacements = [\n  { find: \"68\u00f76=11, 2\", replaceWith: [\"25\u00f77=3, 4\"] },\n  { find: \"65\u00f78=8, 1\", replaceWith: [\"26\u00f76=4, 2\"] },\n  { find: \"62\u00f74=15, 2\", replaceWith: [\"29\u00f78=3, 5\"] },\n  { find: \"99\u00f77=14, 1\", replaceWith: [\"77\u00f78=9, 5\"] },\n  { find: \"20\u00f72=10, 0\", replaceWith: [\"66\u00f74=16, 2\"] },\n  { find: \"97\u00f79=10, 7\", replaceWith: [\"37\u00f75=7, 2\", \"35\u00f74=8, 3\"] },\n  { find: \"23\u00f79=2, 5\", replaceWith: [\"62\u00f72=31, 0\"] },\n  { find: \"98\u00f76=16, 2\", replaceWith: [\"46\u00f72=23, 0\"] },\n  { find: \"16\u00f79=1, 7\", replaceWith: [\"61\u00f73=20, 1\"] },\n  { find: \"66\u00f73=22, 0\", replaceWith: [\"60\u00f77=8, 4\"] },\n  { find: \"48\u00f72=24, 0\", replaceWith: [\"61\u00f79=6, 7\"] },\n  { find: \"66\u00f79=7, 3\", replaceWith: [\"73\u00f77=10, 3\"] },\n  { find: \"69\u00f73=23, 0\", replaceWith: [\"87\u00f73=29, 0\"] },\n  { find: \"42\u00f74=10, 2\", replaceWith: [\"68\u00f72=34, 0\"] },\n  { find: \"37\u00f77=5, 2\", replaceWith: [\"19\u00f76=3, 1\"] },\n  { find: \"99\u00f73=33, 0\", replaceWith: [\"96\u00f76=16, 0\"] },\n  { find: \"74\u00f75=14, 4\", replaceWith: [\"52\u00f76=8, 4\"] },\n  { find: \"65\u00f73=21, 2\", replaceWith: [\"34\u00f77=4, 6\"] },\n  { find: \"24\u00f72=12, 0\", replaceWith: [\"17\u00f78=2, 1\"] },\n  { find: \"29\u00f77=4, 1\", replaceWith: [\"21\u00f72=10, 1\"] },\n  { find: \"46\u00f75=9, 1\", replaceWith: [\"81\u00f77=11, 4\"] },\n  { find: \"56\u00f79=6, 2\", replaceWith: [\"62\u00f76=10, 2\"] },\n  { find: \"61\u00f75=12, 1\", replaceWith: [\"67\u00f75=13, 2\"] },\n  { find: \"41\u00f78=5, 1\", replaceWith: [\"83\u00f75=16, 3\"] },\n];\n\nconst body = context.document.body;\nconst allResults = replacements.map(r =>\n  body.search(r.find, { matchCase: true, matchWholeWord: false })\n);\nallResults.forEach(r => r.load('items'));\nawait context.sync();\n\nreplacements.forEach((r, i) => {\n  const found = allResults[i];\n  for (let j = 0; j < found.items.length; j++) {\n    found.items[j].insertText(r.replaceWith[j], Word.InsertLocation.replace);\n  }\n});\nawait context.sync();\n", "ps1": "# Replace each two-digit-division answer in document order.\n# wdFindContinue=1, wdReplaceOne=1 so duplicate source strings\n# (e.g. \"97\u00f79=10, 7\" appears twice with different targets) are\n# each matched to the correct occurrence in document order.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$null = $find.Execute(\"68\u00f76=11, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"25\u00f77=3, 4\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"65\u00f78=8, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"26\u00f76=4, 2\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"62\u00f74=15, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"29\u00f78=3, 5\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"99\u00f77=14, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"77\u00f78=9, 5\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"20\u00f72=10, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"66\u00f74=16, 2\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"97\u00f79=10, 7\", $false, $false, $false, $false, $false, $true, 1, $false, \"37\u00f75=7, 2\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"23\u00f79=2, 5\", $false, $false, $false, $false, $false, $true, 1, $false, \"62\u00f72=31, 0\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"98\u00f76=16, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"46\u00f72=23, 0\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"16\u00f79=1, 7\", $false, $false, $false, $false, $false, $true, 1, $false, \"61\u00f73=20, 1\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"66\u00f73=22, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"60\u00f77=8, 4\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"48\u00f72=24, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"61\u00f79=6, 7\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"66\u00f79=7, 3\", $false, $false, $false, $false, $false, $true, 1, $false, \"73\u00f77=10, 3\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"69\u00f73=23, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"87\u00f73=29, 0\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"42\u00f74=10, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"68\u00f72=34, 0\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"37\u00f77=5, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"19\u00f76=3, 1\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"97\u00f79=10, 7\", $false, $false, $false, $false, $false, $true, 1, $false, \"35\u00f74=8, 3\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"99\u00f73=33, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"96\u00f76=16, 0\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"74\u00f75=14, 4\", $false, $false, $false, $false, $false, $true, 1, $false, \"52\u00f76=8, 4\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"65\u00f73=21, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"34\u00f77=4, 6\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"24\u00f72=12, 0\", $false, $false, $false, $false, $false, $true, 1, $false, \"17\u00f78=2, 1\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"29\u00f77=4, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"21\u00f72=10, 1\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"46\u00f75=9, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"81\u00f77=11, 4\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"56\u00f79=6, 2\", $false, $false, $false, $false, $false, $true, 1, $false, \"62\u00f76=10, 2\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"61\u00f75=12, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"67\u00f75=13, 2\", 1)\n$find = $d.Content.Find\n$null = $find.Execute(\"41\u00f78=5, 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"83\u00f75=16, 3\", 1)\n\n"}
